$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("VerifyCancelLabelDeferredCredit")
$ws.Range("B2").Value = "Thu Feb 06 14:10:08 IST 2025"
$ws.Range("B3").Value = "Thu Feb 06 14:11:37 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCancelLabelDeferredCorp")
$ws.Range("B2").Value = "Thu Feb 06 14:13:17 IST 2025"
$ws.Range("B3").Value = "Thu Feb 06 14:14:52 IST 2025"
$ws.Range("A3").Value = "Fail"

$ws = $wb.Worksheets.Item("VerifyCancelLabelDeferredPC")
$ws.Range("B2").Value = "Thu Feb 06 14:16:28 IST 2025"
$ws.Range("B3").Value = "Thu Feb 06 14:17:56 IST 2025"
$ws.Range("A3").Value = "Fail"

$ws = $wb.Worksheets.Item("VerifyCancelLabelDeferredPS")
$ws.Range("B2").Value = "Thu Feb 06 14:19:29 IST 2025"
$ws.Range("B3").Value = "Thu Feb 06 14:20:51 IST 2025"
$ws.Range("A3").Value = "Fail"

$ws = $wb.Worksheets.Item("VerifyCreateLabelDeferredCredit")
$ws.Range("B2").Value = "Thu Feb 06 14:22:24 IST 2025"
$ws.Range("B3").Value = "Thu Feb 06 14:23:33 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCreateLabelDeferredCorp")
$ws.Range("B2").Value = "Thu Feb 06 14:24:43 IST 2025"
$ws.Range("B3").Value = "Thu Feb 06 14:25:53 IST 2025"
$ws.Range("A3").Value = "Fail"

$ws = $wb.Worksheets.Item("VerifyCreateLabelDeferredPC")
$ws.Range("B2").Value = "Thu Feb 06 14:27:26 IST 2025"
$ws.Range("B3").Value = "Thu Feb 06 14:28:36 IST 2025"
$ws.Range("A3").Value = "Fail"

$ws = $wb.Worksheets.Item("VerifyCreateLabelDeferredPS")
$ws.Range("B2").Value = "Thu Feb 06 14:30:10 IST 2025"
$ws.Range("B3").Value = "Thu Feb 06 14:31:20 IST 2025"
$ws.Range("A3").Value = "Fail"

$ws = $wb.Worksheets.Item("VerifyEditLabelDeferredCredit")
$ws.Range("B2").Value = "Thu Feb 06 14:32:53 IST 2025"
$ws.Range("B3").Value = "Thu Feb 06 14:34:18 IST 2025"

$ws = $wb.Worksheets.Item("VerifyEditLabelDeferredCorp")
$ws.Range("B2").Value = "Thu Feb 06 14:36:41 IST 2025"
$ws.Range("B3").Value = "Thu Feb 06 14:38:03 IST 2025"
$ws.Range("A3").Value = "Fail"

$ws = $wb.Worksheets.Item("VerifyEditLabelDeferredPC")
$ws.Range("B2").Value = "Thu Feb 06 14:39:35 IST 2025"
$ws.Range("B3").Value = "Thu Feb 06 14:40:58 IST 2025"
$ws.Range("A3").Value = "Fail"

$ws = $wb.Worksheets.Item("VerifyEditLabelDeferredPS")
$ws.Range("B2").Value = "Thu Feb 06 14:42:30 IST 2025"
$ws.Range("B3").Value = "Thu Feb 06 14:43:53 IST 2025"
$ws.Range("A3").Value = "Fail"

$ws = $wb.Worksheets.Item("CCAutoPayPlanCorp")
$ws.Range("B2").Value = "Thu Feb 06 14:45:27 IST 2025"
$ws.Range("B3").Value = "Thu Feb 06 14:46:41 IST 2025"

$ws = $wb.Worksheets.Item("CCAutoPayPlanCredit")
$ws.Range("B2").Value = "Thu Feb 06 14:48:12 IST 2025"
$ws.Range("B3").Value = "Thu Feb 06 14:49:27 IST 2025"

$ws = $wb.Worksheets.Item("CCAutoPayPlanPC")
$ws.Range("B2").Value = "Thu Feb 06 14:50:44 IST 2025"
$ws.Range("B3").Value = "Thu Feb 06 14:51:59 IST 2025"

$ws = $wb.Worksheets.Item("CCAutoPayPlanPS")
$ws.Range("B2").Value = "Thu Feb 06 14:53:31 IST 2025"
$ws.Range("B3").Value = "Thu Feb 06 14:54:45 IST 2025"

$ws = $wb.Worksheets.Item("CCDeferredPlanCorp")
$ws.Range("B2").Value = "Thu Feb 06 14:56:17 IST 2025"
$ws.Range("B3").Value = "Thu Feb 06 14:57:28 IST 2025"

$ws = $wb.Worksheets.Item("CCDeferredPlanCredit")
$ws.Range("B2").Value = "Thu Feb 06 14:59:02 IST 2025"
$ws.Range("B3").Value = "Thu Feb 06 15:00:17 IST 2025"

$ws = $wb.Worksheets.Item("CCDeferredPlanPC")
$ws.Range("B2").Value = "Thu Feb 06 15:01:29 IST 2025"
$ws.Range("B3").Value = "Thu Feb 06 15:02:41 IST 2025"

$ws = $wb.Worksheets.Item("CCDeferredPlanPS")
$ws.Range("B2").Value = "Thu Feb 06 15:04:16 IST 2025"
$ws.Range("B3").Value = "Thu Feb 06 15:05:29 IST 2025"

$ws = $wb.Worksheets.Item("CMCAutoPayPlanCorp")
$ws.Range("B3").Value = "Thu Feb 06 15:08:18 IST 2025"
$ws.Range("A3").Value = "Fail"

$ws = $wb.Worksheets.Item("CMCAutoPayPlanPC")
$ws.Range("B3").Value = "Thu Feb 06 15:13:30 IST 2025"
$ws.Range("A3").Value = "Fail"

$ws = $wb.Worksheets.Item("CMCAutoPayPlanPS")
$ws.Range("B3").Value = "Thu Feb 06 15:16:18 IST 2025"
$ws.Range("A3").Value = "Fail"

$ws = $wb.Worksheets.Item("CMCDeferredPlanCorp")
$ws.Range("B2").Value = "Thu Feb 06 15:17:51 IST 2025"
$ws.Range("B3").Value = "Thu Feb 06 15:19:11 IST 2025"
$ws.Range("A3").Value = "Fail"

$ws = $wb.Worksheets.Item("CMCDeferredPlanCredit")
$ws.Range("B2").Value = "Thu Feb 06 15:20:44 IST 2025"
$ws.Range("B3").Value = "Thu Feb 06 15:22:04 IST 2025"

$ws = $wb.Worksheets.Item("CMCDeferredPlanPC")
$ws.Range("B3").Value = "Thu Feb 06 15:23:25 IST 2025"
$ws.Range("A3").Value = "Fail"

$ws = $wb.Worksheets.Item("CMCDeferredPlanPS")
$ws.Range("B2").Value = "Thu Feb 06 15:24:58 IST 2025"
$ws.Range("B3").Value = "Thu Feb 06 15:26:18 IST 2025"
$ws.Range("A3").Value = "Fail"
